$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 58; existing rows 58-73 shift down to 60-75.
$ws.Rows("58:59").Insert()

# New row 58: Murcott, Primera
$ws.Cells.Item(58, 1).Value = 11
$ws.Cells.Item(58, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(58, 3).Value = "Bíobío"
$ws.Cells.Item(58, 4).Value = 44474
$ws.Cells.Item(58, 5).Value = 8
$ws.Cells.Item(58, 6).Value = "Fruta"
$ws.Cells.Item(58, 7).Value = 100102
$ws.Cells.Item(58, 8).Value = "Cítricos"
$ws.Cells.Item(58, 9).Value = 100102004
$ws.Cells.Item(58, 10).Value = "Mandarina"
$ws.Cells.Item(58, 11).Value = "Murcott"
$ws.Cells.Item(58, 12).Value = "Primera"
$ws.Cells.Item(58, 13).Value = 100
$ws.Cells.Item(58, 14).Value = 5500
$ws.Cells.Item(58, 15).Value = 6000
$ws.Cells.Item(58, 16).Value = 5750
$ws.Cells.Item(58, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(58, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(58, 19).Value = 575
$ws.Cells.Item(58, 20).Value = 10

# New row 59: Murcott, Segunda
$ws.Cells.Item(59, 1).Value = 11
$ws.Cells.Item(59, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(59, 3).Value = "Bíobío"
$ws.Cells.Item(59, 4).Value = 44474
$ws.Cells.Item(59, 5).Value = 8
$ws.Cells.Item(59, 6).Value = "Fruta"
$ws.Cells.Item(59, 7).Value = 100102
$ws.Cells.Item(59, 8).Value = "Cítricos"
$ws.Cells.Item(59, 9).Value = 100102004
$ws.Cells.Item(59, 10).Value = "Mandarina"
$ws.Cells.Item(59, 11).Value = "Murcott"
$ws.Cells.Item(59, 12).Value = "Segunda"
$ws.Cells.Item(59, 13).Value = 50
$ws.Cells.Item(59, 14).Value = 5000
$ws.Cells.Item(59, 15).Value = 5000
$ws.Cells.Item(59, 16).Value = 5000
$ws.Cells.Item(59, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(59, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(59, 19).Value = 500
$ws.Cells.Item(59, 20).Value = 10
